{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of the change (from the OOXML diff):\n//  1. In the \"Update doc gen tools ...\" list item, the run\n//     \" and improve doc gen \" is merged with the word \"process\" (which used\n//     to be a separate run after a _GoBack bookmark) into\n//     \" and improve doc gen process\", followed by the existing \".\" run.\n//     A brand-new list item \"Raise Intel warning level.\" is inserted right\n//     after it (same list style/level) and the _GoBack bookmark now sits at\n//     the end of that new paragraph.\n//  2. The <w:lastRenderedPageBreak/> marker that used to sit on the\n//     \"Fastcall function hooking ...\" run now sits on the preceding\n//     \"Class function hooking ...\" run instead.\n//  3. The <w:lastRenderedPageBreak/> marker that used to sit on the\n//     \"CLR runtime directory support.\" run now sits on the preceding\n//     \"Bound import directory.\" run instead.\n//\n// NOTE: paragraph objects returned by `body.paragraphs` are reseated by\n// index, so after any edit that inserts/removes a paragraph, previously\n// fetched paragraph proxies can silently point at the wrong paragraph.\n// To stay safe we re-locate each target paragraph (by its text) right\n// before editing it, and re-sync between edits.\n\nconst OOXML_WRAPPER =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>__INNER__</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nfunction wrapParagraphs(innerXml) {\n  return OOXML_WRAPPER.replace(\"__INNER__\", innerXml);\n}\n\nasync function findParagraphByText(needle) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n      return paragraphs.items[i];\n    }\n  }\n  throw new Error('Paragraph containing \"' + needle + '\" not found.');\n}\n\n// --- Edit 1: split the \"Update doc gen tools ...\" paragraph into two list\n// items and relocate the _GoBack bookmark to the end of the new one. ---\nconst docGenPara = await findParagraphByText(\"Update doc gen tools\");\nconst docGenReplacement =\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n  '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  \"<w:r><w:t>Update doc gen tools</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> and improve doc gen process</w:t></w:r>' +\n  \"<w:r><w:t>.</w:t></w:r></w:p>\" +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n  '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  \"<w:r><w:t>Raise Intel warning level.</w:t></w:r>\" +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>';\ndocGenPara.insertOoxml(wrapParagraphs(docGenReplacement), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 2: move <w:lastRenderedPageBreak/> from the \"Fastcall function\n// hooking\" run onto the preceding \"Class function hooking\" run. ---\nconst classHookPara = await findParagraphByText(\"Class function hooking\");\nconst classHookReplacement =\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n  '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  \"<w:r><w:lastRenderedPageBreak/><w:t>Class function hooking (ecx preservation).</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> (xchg ecx, [esp]; push ecx)</w:t></w:r></w:p>';\nclassHookPara.insertOoxml(wrapParagraphs(classHookReplacement), Word.InsertLocation.replace);\nawait context.sync();\n\nconst fastcallPara = await findParagraphByText(\"Fastcall function hooking\");\nconst fastcallReplacement =\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n  '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Fastcall function hooking (ecx, edx preservation). </w:t></w:r></w:p>';\nfastcallPara.insertOoxml(wrapParagraphs(fastcallReplacement), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 3: move <w:lastRenderedPageBreak/> from the \"CLR runtime\n// directory support.\" run onto the preceding \"Bound import directory.\" run.\nconst boundImportPara = await findParagraphByText(\"Bound import directory\");\nconst boundImportReplacement =\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/>' +\n  '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  \"<w:r><w:lastRenderedPageBreak/><w:t>Bound import directory.</w:t></w:r></w:p>\";\nboundImportPara.insertOoxml(wrapParagraphs(boundImportReplacement), Word.InsertLocation.replace);\nawait context.sync();\n\nconst clrPara = await findParagraphByText(\"CLR runtime directory support\");\nconst clrReplacement =\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/>' +\n  '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  \"<w:r><w:t>CLR runtime directory support.</w:t></w:r></w:p>\";\nclrPara.insertOoxml(wrapParagraphs(clrReplacement), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Summary of the change (from the OOXML diff):\n#  1. In the \"Update doc gen tools ...\" list item, the run\n#     \" and improve doc gen \" is merged with the word \"process\" (which used\n#     to be a separate run after a _GoBack bookmark) into\n#     \" and improve doc gen process\", followed by the existing \".\" run.\n#     A brand-new list item \"Raise Intel warning level.\" is inserted right\n#     after it (same list style/level) and the _GoBack bookmark now sits at\n#     the end of that new paragraph.\n#  2. The <w:lastRenderedPageBreak/> marker that used to sit on the\n#     \"Fastcall function hooking ...\" run now sits on the preceding\n#     \"Class function hooking ...\" run instead.\n#  3. The <w:lastRenderedPageBreak/> marker that used to sit on the\n#     \"CLR runtime directory support.\" run now sits on the preceding\n#     \"Bound import directory.\" run instead.\n#\n# We use Range.InsertXML (WordprocessingML \"Flat OPC\" fragment) to replace\n# each target paragraph's content precisely, matching the diff exactly.\n# Paragraph indices are re-resolved by searching on text right before each\n# edit (rather than caching stale indices/objects), since the document\n# mutates (paragraph count changes) between edits.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParaIndex([string]$needle) {\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Text -like \"*$needle*\") {\n            return $i\n        }\n    }\n    throw \"Paragraph containing '$needle' not found.\"\n}\n\nfunction New-FlatOpcParagraphXml([string]$innerParagraphsXml) {\n    return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $innerParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# --- Edit 1: split the \"Update doc gen tools ...\" paragraph into two list\n# items and relocate the _GoBack bookmark to the end of the new one. ---\n$docGenIdx = Find-ParaIndex \"Update doc gen tools\"\n$docGenInner = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n    '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>Update doc gen tools</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> and improve doc gen process</w:t></w:r>' +\n    '<w:r><w:t>.</w:t></w:r></w:p>' +\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n    '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>Raise Intel warning level.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n$d.Paragraphs.Item($docGenIdx).Range.InsertXML((New-FlatOpcParagraphXml $docGenInner))\n\n# --- Edit 2: move <w:lastRenderedPageBreak/> from the \"Fastcall function\n# hooking\" run onto the preceding \"Class function hooking\" run. ---\n$classHookIdx = Find-ParaIndex \"Class function hooking\"\n$classHookInner = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n    '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:lastRenderedPageBreak/><w:t>Class function hooking (ecx preservation).</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> (xchg ecx, [esp]; push ecx)</w:t></w:r></w:p>'\n$d.Paragraphs.Item($classHookIdx).Range.InsertXML((New-FlatOpcParagraphXml $classHookInner))\n\n$fastcallIdx = Find-ParaIndex \"Fastcall function hooking\"\n$fastcallInner = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n    '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Fastcall function hooking (ecx, edx preservation). </w:t></w:r></w:p>'\n$d.Paragraphs.Item($fastcallIdx).Range.InsertXML((New-FlatOpcParagraphXml $fastcallInner))\n\n# --- Edit 3: move <w:lastRenderedPageBreak/> from the \"CLR runtime\n# directory support.\" run onto the preceding \"Bound import directory.\" run.\n$boundImportIdx = Find-ParaIndex \"Bound import directory\"\n$boundImportInner = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/>' +\n    '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:lastRenderedPageBreak/><w:t>Bound import directory.</w:t></w:r></w:p>'\n$d.Paragraphs.Item($boundImportIdx).Range.InsertXML((New-FlatOpcParagraphXml $boundImportInner))\n\n$clrIdx = Find-ParaIndex \"CLR runtime directory support\"\n$clrInner = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/>' +\n    '<w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>CLR runtime directory support.</w:t></w:r></w:p>'\n$d.Paragraphs.Item($clrIdx).Range.InsertXML((New-FlatOpcParagraphXml $clrInner))\n"}
